# Slide 25 ("Follow Sets from CPRL: Example 1"), body placeholder:
# split the single run
#   "subprogramDecls = subprogramDecl { subprogramDecl } ."
# into three runs reading
#   "subprogramDecls " + "= { " + "subprogramDecl } ."
# (the grammar rule is corrected to drop the stray leading subprogramDecl).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldText = "subprogramDecls = subprogramDecl { subprogramDecl } ."
$newText = "subprogramDecls = { subprogramDecl } ."

$fullText = $tr.Text
$charIndex = $fullText.IndexOf($oldText)
$runStart = $charIndex + 1

# Rewrite the run's characters in place (keeps the existing run/formatting).
$target = $tr.Characters($runStart, $oldText.Length)
$target.Text = $newText

# Re-split the (now single) run into three runs at the new word
# boundaries, so each piece below becomes its own <a:r>:
#   "subprogramDecls " | "= { " | "subprogramDecl } ."
$part1Len = "subprogramDecls ".Length
$part2Len = "= { ".Length

$part1 = $tr.Characters($runStart, $part1Len)
$part1.Text = $part1.Text

$part2 = $tr.Characters($runStart + $part1Len, $part2Len)
$part2.Text = $part2.Text
